$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8499.75
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 9666.333000000001
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 9666.333000000001
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -10634.333

$ws.Range("H88").Value = 2958
$ws.Range("I88").Value = 1000.6667
$ws.Range("J88").Value = 4426
$ws.Range("K88").Value = 1000.6667
$ws.Range("L88").Value = 4426
$ws.Range("M88").Value = -594.6667
$ws.Range("N88").Value = -5238

$ws.Range("H91").Value = 2958
$ws.Range("I91").Value = 1000.6667
$ws.Range("J91").Value = 4426
$ws.Range("K91").Value = 1000.6667
$ws.Range("L91").Value = 4426
$ws.Range("M91").Value = 403.3333
$ws.Range("N91").Value = -7234

$ws.Range("H129").Value = 833.0599999999999
$ws.Range("I129").Value = 343.83334
$ws.Range("J129").Value = 899.7727
$ws.Range("K129").Value = 1031.50002
$ws.Range("L129").Value = 2699.3181
$ws.Range("M129").Value = 3968.49998
$ws.Range("N129").Value = -12699.3181

$ws.Range("H133").Value = 51545
$ws.Range("J133").Value = 51545
$ws.Range("L133").Value = 51545
$ws.Range("N133").Value = -61665

$ws.Range("H136").Value = 46611
$ws.Range("J136").Value = 46611
$ws.Range("L136").Value = 46611
$ws.Range("N136").Value = -56811

$ws.Range("H137").Value = 3579.923
$ws.Range("I137").Value = 2707
$ws.Range("K137").Value = 8121
$ws.Range("M137").Value = -5571

$ws.Range("H139").Value = 33180
$ws.Range("J139").Value = 33180
$ws.Range("L139").Value = 33180
$ws.Range("N139").Value = -43460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1313.9
$ws.Range("I61").Value = 1381.762
$ws.Range("J61").Value = 1155.5555
$ws.Range("K61").Value = 1381.762
$ws.Range("L61").Value = 1155.5555
$ws.Range("M61").Value = -1169.762
$ws.Range("N61").Value = -1579.5555

$ws.Range("H74").Value = 4959.25
$ws.Range("I74").Value = 6461.231
$ws.Range("J74").Value = 3184.182
$ws.Range("K74").Value = 6461.231
$ws.Range("L74").Value = 3184.182
$ws.Range("M74").Value = -5587.231
$ws.Range("N74").Value = -4932.182

$ws.Range("H77").Value = 4959.25
$ws.Range("I77").Value = 6461.231
$ws.Range("J77").Value = 3184.182
$ws.Range("K77").Value = 32306.155
$ws.Range("L77").Value = 15920.91
$ws.Range("M77").Value = -27938.155
$ws.Range("N77").Value = -24656.91

$ws.Range("H102").Value = 2353.818
$ws.Range("I102").Value = 1973.6666
$ws.Range("J102").Value = 2810
$ws.Range("K102").Value = 1973.6666
$ws.Range("L102").Value = 2810
$ws.Range("M102").Value = -351.6666
$ws.Range("N102").Value = -6054

$ws.Range("H132").Value = 3272.6191
$ws.Range("I132").Value = 1988.8889
$ws.Range("J132").Value = 4235.4165
$ws.Range("K132").Value = 5966.6667
$ws.Range("L132").Value = 12706.2495
$ws.Range("M132").Value = -3436.6667
$ws.Range("N132").Value = -17766.2495

$ws.Range("H136").Value = 1313.9
$ws.Range("I136").Value = 1381.762
$ws.Range("J136").Value = 1155.5555
$ws.Range("K136").Value = 4145.286
$ws.Range("L136").Value = 3466.6665
$ws.Range("M136").Value = -1595.286
$ws.Range("N136").Value = -8566.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 3462.6667
$ws.Range("I26").Value = 3462.6667
$ws.Range("K26").Value = 3462.6667
$ws.Range("M26").Value = -3170.6667

$ws.Range("H96").Value = 4280
$ws.Range("I96").Value = 4280
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4280
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1534
$ws.Range("N96").Value = $null

$ws.Range("H107").Value = 2262.7144
$ws.Range("I107").Value = 2139.8333
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2139.8333
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -219.8332999999998
$ws.Range("N107").Value = -6840

$ws.Range("H122").Value = 41780
$ws.Range("J122").Value = 41780
$ws.Range("L122").Value = 41780
$ws.Range("N122").Value = -51580

$ws.Range("H126").Value = 41868.57
$ws.Range("J126").Value = 41868.57
$ws.Range("L126").Value = 41868.57
$ws.Range("N126").Value = -51748.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 43653.11
$ws.Range("J20").Value = 43653.11
$ws.Range("L20").Value = 43653.11
$ws.Range("N20").Value = -44125.11

$ws.Range("H30").Value = 43653.11
$ws.Range("J30").Value = 43653.11
$ws.Range("L30").Value = 43653.11
$ws.Range("N30").Value = -43835.11

$ws.Range("H31").Value = 3446.3125
$ws.Range("I31").Value = 1194.6364
$ws.Range("J31").Value = 8400
$ws.Range("K31").Value = 1194.6364
$ws.Range("L31").Value = 8400
$ws.Range("M31").Value = -899.6364000000001
$ws.Range("N31").Value = -8990

$ws.Range("H34").Value = 3446.3125
$ws.Range("I34").Value = 1194.6364
$ws.Range("J34").Value = 8400
$ws.Range("K34").Value = 1194.6364
$ws.Range("L34").Value = 8400
$ws.Range("M34").Value = -992.6364000000001
$ws.Range("N34").Value = -8804

$ws.Range("H99").Value = 10004044
$ws.Range("I99").Value = 33335230
$ws.Range("J99").Value = 4964.5713
$ws.Range("K99").Value = 33335230
$ws.Range("L99").Value = 4964.5713
$ws.Range("M99").Value = -33333732
$ws.Range("N99").Value = -7960.5713

$ws.Range("H105").Value = 1741.3
$ws.Range("I105").Value = 1345.5
$ws.Range("K105").Value = 1345.5
$ws.Range("M105").Value = 401.5

$ws.Range("H126").Value = 10004044
$ws.Range("I126").Value = 33335230
$ws.Range("J126").Value = 4964.5713
$ws.Range("K126").Value = 100005690
$ws.Range("L126").Value = 14893.7139
$ws.Range("M126").Value = -100003220
$ws.Range("N126").Value = -19833.7139

$ws.Range("H128").Value = 43653.11
$ws.Range("J128").Value = 43653.11
$ws.Range("L128").Value = 43653.11
$ws.Range("N128").Value = -53613.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 306.6
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 314.14285
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 942.4285500000001
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -3438.42855

$ws.Range("H113").Value = 666.1702
$ws.Range("I113").Value = 662.0294
$ws.Range("J113").Value = 677
$ws.Range("K113").Value = 1986.0882
$ws.Range("L113").Value = 2031
$ws.Range("M113").Value = 183.9117999999999
$ws.Range("N113").Value = -6371

$ws.Range("H121").Value = 2586.7222
$ws.Range("I121").Value = 230
$ws.Range("J121").Value = 2631.1887
$ws.Range("K121").Value = 690
$ws.Range("L121").Value = 7893.5661
$ws.Range("M121").Value = 620
$ws.Range("N121").Value = -10513.5661

$ws.Range("H131").Value = 5883131.5
$ws.Range("I131").Value = 125000424
$ws.Range("J131").Value = 795.9752999999999
$ws.Range("K131").Value = 375001272
$ws.Range("L131").Value = 2387.9259
$ws.Range("M131").Value = -374996232
$ws.Range("N131").Value = -12467.9259

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25002732
$ws.Range("I80").Value = 35716890
$ws.Range("J80").Value = 3033.3333
$ws.Range("K80").Value = 35716890
$ws.Range("L80").Value = 3033.3333
$ws.Range("M80").Value = -35715892
$ws.Range("N80").Value = -5029.3333

$ws.Range("H83").Value = 25002732
$ws.Range("I83").Value = 35716890
$ws.Range("J83").Value = 3033.3333
$ws.Range("K83").Value = 178584450
$ws.Range("L83").Value = 15166.6665
$ws.Range("M83").Value = -178579458
$ws.Range("N83").Value = -25150.6665

$ws.Range("H102").Value = 2421.3125
$ws.Range("I102").Value = 1587
$ws.Range("J102").Value = 6036.6665
$ws.Range("K102").Value = 1587
$ws.Range("L102").Value = 6036.6665
$ws.Range("M102").Value = 35
$ws.Range("N102").Value = -9280.666499999999

$ws.Range("H113").Value = 1850
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340

$ws.Range("H126").Value = 3237.07
$ws.Range("I126").Value = 2901.274
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8703.822
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6233.822
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 4377.6
$ws.Range("I132").Value = 2682.5715
$ws.Range("K132").Value = 8047.7145
$ws.Range("M132").Value = -5517.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996

$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984

$ws.Range("H122").Value = 2855.1538
$ws.Range("I122").Value = 1460.5
$ws.Range("K122").Value = 4381.5
$ws.Range("M122").Value = -1931.5

$ws.Range("H132").Value = 5382.72
$ws.Range("I132").Value = 1797.6923
$ws.Range("J132").Value = 9266.5
$ws.Range("K132").Value = 5393.0769
$ws.Range("L132").Value = 27799.5
$ws.Range("M132").Value = -2863.0769
$ws.Range("N132").Value = -32859.5

$ws.Range("H133").Value = 33661.25
$ws.Range("J133").Value = 33661.25
$ws.Range("L133").Value = 33661.25
$ws.Range("N133").Value = -38721.25

$ws.Range("H136").Value = 3857.7856
$ws.Range("I136").Value = 1179.8572
$ws.Range("J136").Value = 6535.7144
$ws.Range("K136").Value = 3539.5716
$ws.Range("L136").Value = 19607.1432
$ws.Range("M136").Value = -989.5715999999998
$ws.Range("N136").Value = -24707.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 25335.334
$ws.Range("I18").Value = 3000
$ws.Range("J18").Value = 36503
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 36503
$ws.Range("M18").Value = -2827
$ws.Range("N18").Value = -36849

$ws.Range("H31").Value = 26000
$ws.Range("J31").Value = 26000
$ws.Range("L31").Value = 26000
$ws.Range("N31").Value = -26696

$ws.Range("H132").Value = 10418468
$ws.Range("I132").Value = 1090.2
$ws.Range("J132").Value = 47623390
$ws.Range("K132").Value = 3270.6
$ws.Range("L132").Value = 142870170
$ws.Range("M132").Value = -740.6000000000004
$ws.Range("N132").Value = -142875230

$ws.Range("H136").Value = 5556.857
$ws.Range("I136").Value = 3639.5334
$ws.Range("J136").Value = 10350.167
$ws.Range("K136").Value = 10918.6002
$ws.Range("L136").Value = 31050.501
$ws.Range("M136").Value = -8368.600199999999
$ws.Range("N136").Value = -36150.501
